$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" date column (C) for rows 2-9 from 45170 (2023-09-01) to 45174 (2023-09-05)
$ws.Range("C2:C9").Value = 45174
